$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (bold)
$ws.Range("A2").Value = "SR.No"
$ws.Range("B2").Value = "Name"
$ws.Range("A2:B2").Font.Bold = $true

# Data rows
$ws.Range("A4").Value = 1
$ws.Range("B4").Value = "Nitin Singh"
$ws.Range("A5").Value = 2
$ws.Range("B5").Value = "Pravin Shukla"
$ws.Range("A6").Value = 3
$ws.Range("B6").Value = "Maddy Singh"

# Page setup - portrait orientation
$ws.PageSetup.Orientation = 1

# Move the active selection to D11, as in the saved workbook
$ws.Range("D11").Select() | Out-Null
